$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1 header: average_doctor / average_doctor_old swapped position
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Data rows 4-13: recalculated statistic values (harvard case classification)

# Row 4
$ws.Range("AI4").Value = 0.21
$ws.Range("AJ4").Value = 0.066
$ws.Range("AK4").Value = 0.257
$ws.Range("AU4").Value = 0.143
$ws.Range("AV4").Value = 0.027
$ws.Range("AW4").Value = 0.165
$ws.Range("BA4").Value = 1.938
$ws.Range("BB4").Value = 0.169
$ws.Range("BC4").Value = 0.411
$ws.Range("BG4").Value = 0.723
$ws.Range("BH4").Value = 0.147
$ws.Range("BI4").Value = 0.383
$ws.Range("BM4").Value = 0.676
$ws.Range("BN4").Value = 0.09
$ws.Range("BO4").Value = 0.3
$ws.Range("BP4").Value = 0.646
$ws.Range("BQ4").Value = 0.656
$ws.Range("E4").Value = 0.384
$ws.Range("F4").Value = 0.076
$ws.Range("G4").Value = 0.275
$ws.Range("N4").Value = 0.385
$ws.Range("O4").Value = 0.061
$ws.Range("P4").Value = 0.247
$ws.Range("W4").Value = 0.22
$ws.Range("X4").Value = 0.104
$ws.Range("Y4").Value = 0.323

# Row 5
$ws.Range("AI5").Value = 0.25
$ws.Range("AJ5").Value = 0.095
$ws.Range("AK5").Value = 0.309
$ws.Range("AU5").Value = 0.281
$ws.Range("AV5").Value = 0.091
$ws.Range("AW5").Value = 0.302
$ws.Range("BA5").Value = 1.337
$ws.Range("BB5").Value = 0.08
$ws.Range("BC5").Value = 0.282
$ws.Range("BG5").Value = 0.393
$ws.Range("BH5").Value = 0.048
$ws.Range("BI5").Value = 0.218
$ws.Range("BM5").Value = 0.5679999999999999
$ws.Range("BN5").Value = 0.079
$ws.Range("BO5").Value = 0.28
$ws.Range("BP5").Value = 0.446
$ws.Range("BQ5").Value = 0.452
$ws.Range("E5").Value = 0.496
$ws.Range("F5").Value = 0.095
$ws.Range("G5").Value = 0.309
$ws.Range("N5").Value = 0.747
$ws.Range("O5").Value = 0.08699999999999999
$ws.Range("P5").Value = 0.296
$ws.Range("W5").Value = 0.229
$ws.Range("X5").Value = 0.117
$ws.Range("Y5").Value = 0.342

# Row 6
$ws.Range("AI6").Value = 0.228
$ws.Range("AU6").Value = 0.19
$ws.Range("BA6").Value = 1.569
$ws.Range("BG6").Value = 0.509
$ws.Range("BM6").Value = 0.617
$ws.Range("BP6").Value = 0.523
$ws.Range("BQ6").Value = 0.532
$ws.Range("E6").Value = 0.433
$ws.Range("N6").Value = 0.508
$ws.Range("W6").Value = 0.224

# Row 7
$ws.Range("AI7").Value = 0.241
$ws.Range("AU7").Value = 0.236
$ws.Range("BA7").Value = 1.419
$ws.Range("BG7").Value = 0.432
$ws.Range("BM7").Value = 0.587
$ws.Range("BP7").Value = 0.473
$ws.Range("BQ7").Value = 0.48
$ws.Range("E7").Value = 0.469
$ws.Range("N7").Value = 0.629
$ws.Range("W7").Value = 0.227

# Row 8
$ws.Range("AI8").Value = 0.229
$ws.Range("AJ8").Value = 0.095
$ws.Range("AK8").Value = 0.308
$ws.Range("AU8").Value = 0.225
$ws.Range("AV8").Value = 0.07199999999999999
$ws.Range("AW8").Value = 0.268
$ws.Range("BA8").Value = 1.721
$ws.Range("BB8").Value = 0.139
$ws.Range("BC8").Value = 0.373
$ws.Range("BG8").Value = 0.549
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.327
$ws.Range("BM8").Value = 0.703
$ws.Range("BN8").Value = 0.074
$ws.Range("BO8").Value = 0.272
$ws.Range("BP8").Value = 0.574
$ws.Range("BQ8").Value = 0.584
$ws.Range("E8").Value = 0.528
$ws.Range("F8").Value = 0.126
$ws.Range("G8").Value = 0.354
$ws.Range("N8").Value = 0.739
$ws.Range("O8").Value = 0.075
$ws.Range("P8").Value = 0.275
$ws.Range("W8").Value = 0.228
$ws.Range("X8").Value = 0.11
$ws.Range("Y8").Value = 0.332

# Row 9
$ws.Range("AI9").Value = 0.128
$ws.Range("AJ9").Value = 0.111
$ws.Range("AK9").Value = 0.334
$ws.Range("BA9").Value = 1.68
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.574
$ws.Range("BH9").Value = 0.244
$ws.Range("BI9").Value = 0.494
$ws.Range("BM9").Value = 0.638
$ws.Range("BN9").Value = 0.231
$ws.Range("BO9").Value = 0.48
$ws.Range("BP9").Value = 0.5600000000000001
$ws.Range("BQ9").Value = 0.5570000000000001
$ws.Range("E9").Value = 0.447
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("N9").Value = 0.596
$ws.Range("O9").Value = 0.241
$ws.Range("P9").Value = 0.491
$ws.Range("W9").Value = 0.128
$ws.Range("X9").Value = 0.111
$ws.Range("Y9").Value = 0.334

# Row 10
$ws.Range("AI10").Value = 0.255
$ws.Range("AJ10").Value = 0.19
$ws.Range("AK10").Value = 0.436
$ws.Range("AU10").Value = 0.213
$ws.Range("AV10").Value = 0.167
$ws.Range("AW10").Value = 0.409
$ws.Range("BA10").Value = 1.979
$ws.Range("BG10").Value = 0.617
$ws.Range("BH10").Value = 0.236
$ws.Range("BI10").Value = 0.486
$ws.Range("BM10").Value = 0.851
$ws.Range("BN10").Value = 0.127
$ws.Range("BO10").Value = 0.356
$ws.Range("BP10").Value = 0.66
$ws.Range("BQ10").Value = 0.6840000000000001
$ws.Range("E10").Value = 0.574
$ws.Range("F10").Value = 0.244
$ws.Range("G10").Value = 0.494
$ws.Range("N10").Value = 0.8090000000000001
$ws.Range("O10").Value = 0.155
$ws.Range("P10").Value = 0.393
$ws.Range("W10").Value = 0.277
$ws.Range("X10").Value = 0.2
$ws.Range("Y10").Value = 0.447

# Row 11
$ws.Range("AI11").Value = 0.255
$ws.Range("AJ11").Value = 0.19
$ws.Range("AK11").Value = 0.436
$ws.Range("AU11").Value = 0.319
$ws.Range("AV11").Value = 0.217
$ws.Range("AW11").Value = 0.466
$ws.Range("BA11").Value = 1.979
$ws.Range("BG11").Value = 0.617
$ws.Range("BH11").Value = 0.236
$ws.Range("BI11").Value = 0.486
$ws.Range("BM11").Value = 0.851
$ws.Range("BN11").Value = 0.127
$ws.Range("BO11").Value = 0.356
$ws.Range("BP11").Value = 0.66
$ws.Range("BQ11").Value = 0.6840000000000001
$ws.Range("E11").Value = 0.596
$ws.Range("F11").Value = 0.241
$ws.Range("G11").Value = 0.491
$ws.Range("N11").Value = 0.851
$ws.Range("O11").Value = 0.127
$ws.Range("P11").Value = 0.356
$ws.Range("W11").Value = 0.277
$ws.Range("X11").Value = 0.2
$ws.Range("Y11").Value = 0.447

# Row 12
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.667
$ws.Range("AV12").Value = 1.689
$ws.Range("AW12").Value = 1.3
$ws.Range("BA12").Value = 3.519
$ws.Range("BB12").Value = 0.193
$ws.Range("BC12").Value = 0.439
$ws.Range("BG12").Value = 1.069
$ws.Range("BH12").Value = 0.064
$ws.Range("BI12").Value = 0.253
$ws.Range("BM12").Value = 1.325
$ws.Range("BN12").Value = 0.369
$ws.Range("BO12").Value = 0.608
$ws.Range("BP12").Value = 1.173
$ws.Range("BQ12").Value = 1.23
$ws.Range("E12").Value = 1.429
$ws.Range("F12").Value = 0.673
$ws.Range("G12").Value = 0.821
$ws.Range("N12").Value = 1.714
$ws.Range("O12").Value = 1.728
$ws.Range("P12").Value = 1.314
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863

# Row 13
$ws.Range("AI13").Value = 1.352
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK13").Value = 0.633
$ws.Range("AU13").Value = 2.437
$ws.Range("AV13").Value = 1.363
$ws.Range("AW13").Value = 1.167
$ws.Range("BA13").Value = 2.456
$ws.Range("BB13").Value = 0.303
$ws.Range("BC13").Value = 0.551
$ws.Range("BG13").Value = 0.597
$ws.Range("BH13").Value = 0.051
$ws.Range("BI13").Value = 0.225
$ws.Range("BM13").Value = 0.977
$ws.Range("BN13").Value = 0.29
$ws.Range("BO13").Value = 0.538
$ws.Range("BP13").Value = 0.819
$ws.Range("BQ13").Value = 0.782
$ws.Range("E13").Value = 1.657
$ws.Range("F13").Value = 0.711
$ws.Range("G13").Value = 0.843
$ws.Range("N13").Value = 2.381
$ws.Range("O13").Value = 1.153
$ws.Range("P13").Value = 1.074
$ws.Range("W13").Value = 1.055
$ws.Range("X13").Value = 0.172
$ws.Range("Y13").Value = 0.415
